$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing space from the student's surname "KOH " -> "KOH"
$ws.Range("B3").Value = "KOH"

# Move selection to B4 (as would naturally happen after editing B3 and pressing Enter)
$ws.Range("B4").Select()
